$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 28: update title and link
$ws.Range("D28").Value = "Intel의 Realsense 대신 사용할 3D Sensor들은 뭐가 있을까?"
$ws.Range("E28").Value = "https://ropiens.tistory.com/147"

# Row 36: update title and link
$ws.Range("D36").Value = "Active Learning in Semiconductor Manufacturing"
$ws.Range("E36").Value = "http://dmqm.korea.ac.kr/activity/seminar/331"

# Row 37: update title and link
$ws.Range("D37").Value = "[Paper Review] PUMAD : PU Metric Learning for Anomaly Detection"
$ws.Range("E37").Value = "http://dsba.korea.ac.kr/seminar/?uid=1813&mod=document&pageid=1"
